# "can save and load shopping lists"
# Adds a second worksheet ("List2") to the workbook, positioned right
# after the existing "Data" sheet, representing a saved shopping list.
# The new sheet becomes the active/selected sheet (as it would be right
# after the app created & switched to it), with two entries ("Steak")
# in column A and the view's cursor left at E28.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("Data")

# Insert the new sheet immediately after "Data"
$listSheet = $wb.Worksheets.Add($null, $dataSheet)
$listSheet.Name = "List2"

# Shopping-list contents
$listSheet.Range("A1").Value = "Steak"
$listSheet.Range("A2").Value = "Steak"

# Match the outline defaults used for a freshly created sheet
$listSheet.Outline.SummaryRow = 1
$listSheet.Outline.SummaryColumn = 1

# Make the new list sheet the active one, with the same cursor position
# Excel leaves behind after the list was written
$listSheet.Activate() | Out-Null
$listSheet.Range("E28").Select() | Out-Null
